$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Profit Margin" header label in C10, matching the style of B10 ("Total Revenue")
$ws.Range("C10").Value = "Profit Margin"
$ws.Range("B10").Copy()
$ws.Range("C10").PasteSpecial(-4122)  # xlPasteFormats

# Add the Profit Margin total formula in C11, matching the style/number format of B11
$ws.Range("C11").Formula = "=SUMPRODUCT(B3:B7,C3:C7,D3:D7)"
$ws.Range("B11").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats

# Match style of C12 to the currency format column used in row 11 (empty placeholder row)
$ws.Range("C12").Value = $null
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(12).RowHeight = 15.75

$excel.CutCopyMode = $false
$ws.Range("D10").Select()
